# Update the "Datos" ratios table with refreshed Yahoo Finance scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - AAPL
$ws.Range("B2").Value = 1.23
$ws.Range("E2").Value = 2164000000000
$ws.Range("F2").Value = 38.38
$ws.Range("G2").Value = 430.06

# Row 3 - AMD
$ws.Range("B3").Value = 2.29
$ws.Range("C3").Value = 76.15
$ws.Range("E3").Value = 100993000000
$ws.Range("F3").Value = 167.35

# Row 4 - BABA
$ws.Range("B4").Value = 1.6
$ws.Range("C4").Value = 301.35
$ws.Range("E4").Value = 767145000000
$ws.Range("F4").Value = 83.49

# Row 5 - BBD
$ws.Range("B5").Value = 1.29
$ws.Range("E5").Value = 30903000000
$ws.Range("F5").Value = 5.24

# Row 6 - CVX
$ws.Range("B6").Value = 1.26
$ws.Range("E6").Value = 158309000000
$ws.Range("G6").Value = 100.27

# Row 7 - FB
$ws.Range("B7").Value = 1.21
$ws.Range("C7").Value = 284.24
$ws.Range("E7").Value = 865785000000
$ws.Range("F7").Value = 37.16

# Row 8 - GGAL
$ws.Range("B8").Value = 2.06
$ws.Range("C8").Value = 12.11
$ws.Range("E8").Value = 2558000000

# Row 9 - GLOB
$ws.Range("B9").Value = 1.39
$ws.Range("C9").Value = 166.67
$ws.Range("D9").Value = 1.35
$ws.Range("E9").Value = 7078000000
$ws.Range("F9").Value = 132.63

# Row 10 - JPM
$ws.Range("B10").Value = 1.1599999999999999
$ws.Range("E10").Value = 301987000000
$ws.Range("F10").Value = 13.27
$ws.Range("G10").Value = 115.04

# Row 11 - MELI
$ws.Range("B11").Value = 1.71
$ws.Range("C11").Value = 1170.95
$ws.Range("D11").Value = -3.4
$ws.Range("E11").Value = 60444000000

# Row 12 - MSFT
$ws.Range("B12").Value = 0.9
$ws.Range("E12").Value = 1674000000000
$ws.Range("F12").Value = 38.39

# Row 13 - NFLX
$ws.Range("B13").Value = 0.94
$ws.Range("C13").Value = 515.36
$ws.Range("E13").Value = 241469000000
$ws.Range("F13").Value = 92.46

# Row 14 - NIO
$ws.Range("C14").Value = 13.02
$ws.Range("E14").Value = 2421000000

# Row 15 - NKE
$ws.Range("B15").Value = 0.78
$ws.Range("E15").Value = 173975000000
$ws.Range("F15").Value = 69.709999999999994
$ws.Range("G15").Value = 112.79

# Row 16 - NVDA
$ws.Range("B16").Value = 1.35
$ws.Range("D16").Value = 5.44
$ws.Range("E16").Value = 315238000000
$ws.Range("F16").Value = 93.85
$ws.Range("G16").Value = 520.92999999999995

# Row 17 - PAM
$ws.Range("B17").Value = 1.26
$ws.Range("C17").Value = 14.76
$ws.Range("E17").Value = 1241000000

# Row 18 - PBR
$ws.Range("B18").Value = 1.92
$ws.Range("E18").Value = 51537000000
$ws.Range("G18").Value = 13.75

# Row 19 - PFE
$ws.Range("B19").Value = 0.73
$ws.Range("E19").Value = 211439000000
$ws.Range("F19").Value = 15.06
$ws.Range("G19").Value = 41.79

# Row 20 - ROKU
$ws.Range("B20").Value = 1.89
$ws.Range("C20").Value = 163.9
$ws.Range("D20").Value = -1.1599999999999999
$ws.Range("E20").Value = 20487000000

# Row 21 - SHOP
$ws.Range("B21").Value = 1.58
$ws.Range("C21").Value = 1094.83
$ws.Range("E21").Value = 130336000000

# Row 22 - TSLA
$ws.Range("B22").Value = 1.3
$ws.Range("C22").Value = 1340.18
$ws.Range("E22").Value = 401269000000
$ws.Range("F22").Value = 1114.48

# Row 23 - XOM
$ws.Range("B23").Value = 1.3
$ws.Range("E23").Value = 169171000000
$ws.Range("F23").Value = 23.83
$ws.Range("G23").Value = 48

# Row 24 - YPF
$ws.Range("B24").Value = 1.75
$ws.Range("E24").Value = 2187000000
$ws.Range("F24").Value = 1.78
$ws.Range("G24").Value = 7.13

# Restore the active selection to E23 (single cell) as in the final file.
$ws.Range("E23").Select()
